$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1902.3334
$ws.Range("J17").Value = 1902.3334
$ws.Range("L17").Value = 5707.0002
$ws.Range("N17").Value = -6043.0002

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1128.2222
$ws.Range("I41").Value = 372
$ws.Range("J41").Value = 3775
$ws.Range("K41").Value = 372
$ws.Range("L41").Value = 3775
$ws.Range("M41").Value = 68
$ws.Range("N41").Value = -4655

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7680.857
$ws.Range("I62").Value = 3960
$ws.Range("J62").Value = 8301
$ws.Range("K62").Value = 3960
$ws.Range("L62").Value = 8301
$ws.Range("M62").Value = -3336
$ws.Range("N62").Value = -9549

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7680.857
$ws.Range("I65").Value = 3960
$ws.Range("J65").Value = 8301
$ws.Range("K65").Value = 19800
$ws.Range("L65").Value = 41505
$ws.Range("M65").Value = -16680
$ws.Range("N65").Value = -47745

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 99999
$ws.Range("J87").Value = 99999
$ws.Range("L87").Value = 99999
$ws.Range("N87").Value = -102495

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 99999
$ws.Range("J90").Value = 99999
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477

# ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 7079.8887
$ws.Range("J99").Value = 7358.625
$ws.Range("L99").Value = 22075.875
$ws.Range("N99").Value = -25071.875

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8500
$ws.Range("I113").Value = 8501.25
$ws.Range("J113").Value = 8497.5
$ws.Range("K113").Value = 8501.25
$ws.Range("L113").Value = 8497.5
$ws.Range("M113").Value = -5247.25
$ws.Range("N113").Value = -15005.5

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1030.4445
$ws.Range("I2").Value = 1030.4445
$ws.Range("K2").Value = 1030.4445
$ws.Range("M2").Value = -917.4445000000001

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1288.5264
$ws.Range("I32").Value = 1155.375
$ws.Range("J32").Value = 1998.6666
$ws.Range("K32").Value = 1155.375
$ws.Range("L32").Value = 1998.6666
$ws.Range("M32").Value = -868.375
$ws.Range("N32").Value = -2572.6666

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 789.5
$ws.Range("I97").Value = 773.3077
$ws.Range("K97").Value = 773.3077
$ws.Range("M97").Value = -277.3077

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1030.4445
$ws.Range("I116").Value = 1030.4445
$ws.Range("K116").Value = 1030.4445
$ws.Range("M116").Value = 1263.5555

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1030.4445
$ws.Range("I3").Value = 1030.4445
$ws.Range("K3").Value = 1030.4445
$ws.Range("M3").Value = -916.4445000000001

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4753.4614
$ws.Range("I86").Value = 1399.2858
$ws.Range("K86").Value = 1399.2858
$ws.Range("M86").Value = -276.2858000000001

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4753.4614
$ws.Range("I89").Value = 1399.2858
$ws.Range("K89").Value = 6996.429
$ws.Range("M89").Value = -1380.429

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1268.0714
$ws.Range("I99").Value = 1057.9231
$ws.Range("K99").Value = 1057.9231
$ws.Range("M99").Value = 440.0769

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1685.4
$ws.Range("J105").Value = 1969.2
$ws.Range("L105").Value = 1969.2
$ws.Range("N105").Value = -5463.2

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 138130.62
$ws.Range("I16").Value = 157756.42
$ws.Range("K16").Value = 157756.42
$ws.Range("M16").Value = -157469.42

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6168.5293
$ws.Range("I31").Value = 1858.125
$ws.Range("K31").Value = 1858.125
$ws.Range("M31").Value = -1563.125

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6168.5293
$ws.Range("I34").Value = 1858.125
$ws.Range("K34").Value = 1858.125
$ws.Range("M34").Value = -1656.125

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2083.1667
$ws.Range("I99").Value = 2090.7273
$ws.Range("K99").Value = 2090.7273
$ws.Range("M99").Value = -592.7273

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 138130.62
$ws.Range("I113").Value = 157756.42
$ws.Range("K113").Value = 157756.42
$ws.Range("M113").Value = -155586.42

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2083.1667
$ws.Range("I126").Value = 2090.7273
$ws.Range("K126").Value = 6272.1819
$ws.Range("M126").Value = -3802.1819

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7750
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 15000
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 45000
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -46622

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 7750
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 15000
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 135000
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -143112

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4269.625
$ws.Range("I80").Value = 3888.8333
$ws.Range("K80").Value = 11666.4999
$ws.Range("M80").Value = -10730.4999

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4269.625
$ws.Range("I83").Value = 3888.8333
$ws.Range("K83").Value = 34999.4997
$ws.Range("M83").Value = -30319.4997

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 7980
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 560.5625
$ws.Range("I97").Value = 379
$ws.Range("J97").Value = 960
$ws.Range("K97").Value = 379
$ws.Range("L97").Value = 960
$ws.Range("M97").Value = 117
$ws.Range("N97").Value = -1952

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1955.75
$ws.Range("I102").Value = 1736.2142
$ws.Range("K102").Value = 1736.2142
$ws.Range("M102").Value = -114.2141999999999

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8060

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5258.875
$ws.Range("I40").Value = 4581
$ws.Range("J40").Value = 10004
$ws.Range("K40").Value = 4581
$ws.Range("L40").Value = 10004
$ws.Range("M40").Value = -4445
$ws.Range("N40").Value = -10276

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5116.5
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5116.5
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1001
$ws.Range("I107").Value = 1001
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3003
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1083
$ws.Range("N107").ClearContents()

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4121.242
$ws.Range("I126").Value = 2350.3
$ws.Range("J126").Value = 6845.769
$ws.Range("K126").Value = 7050.900000000001
$ws.Range("L126").Value = 20537.307
$ws.Range("M126").Value = -4580.900000000001
$ws.Range("N126").Value = -25477.307

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1585.2222
$ws.Range("I132").Value = 1585.2222
$ws.Range("K132").Value = 4755.6666
$ws.Range("M132").Value = -2225.6666
